$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest achievable quantized values given runtime 1/6 rounding)
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Update cell values
$ws.Range("A1").Value = -0.12831060389814297
$ws.Range("B1").Value = 0.12824627590856608
$ws.Range("A2").Value = -0.10613994171847363
$ws.Range("B2").Value = 0.10594806856791994
$ws.Range("A3").Value = -0.062441373547418522
$ws.Range("B3").Value = 0.062318912095170376
$ws.Range("A4").Value = -0.054318912112252704
$ws.Range("B4").Value = 0.053783936798721399
$ws.Range("A5").Value = -0.050783936807293983
$ws.Range("B5").Value = 0.048955049268253781
$ws.Range("A6").Value = -0.041835639482542319
$ws.Range("B6").Value = 0.041265424208017265
$ws.Range("A7").Value = -0.031265424231284644
$ws.Range("B7").Value = 0.03112230113408021
$ws.Range("A8").Value = -0.021122301158092771
$ws.Range("B8").Value = 0.020855089219321332
$ws.Range("A9").Value = -0.018855089229491639
$ws.Range("B9").Value = 0.018631840818514878
$ws.Range("A10").Value = -0.027412126027615002
$ws.Range("B10").Value = 0.027395471896751289
$ws.Range("A11").Value = -0.024395471909585353
$ws.Range("B11").Value = 0.024368547108363892
$ws.Range("A12").Value = -0.020868547122380843
$ws.Range("B12").Value = 0.020672939597704776
$ws.Range("A13").Value = -0.01717293961296118
$ws.Range("B13").Value = 0.01708359156231154
$ws.Range("A14").Value = -0.0090835915861307015
$ws.Range("B14").Value = 0.0090541547483837803
$ws.Range("A15").Value = -0.0080541547596091334
$ws.Range("B15").Value = 0.0080351248432730671
$ws.Range("A16").Value = -0.0060351248565284621
$ws.Range("B16").Value = 0.0060037631217979559
$ws.Range("A17").Value = -0.0040037631353513348
$ws.Range("B17").Value = 0.0039999999827058375
$ws.Range("A18").Value = -0.016106334201307249
$ws.Range("B18").Value = 0.016091856038233487
$ws.Range("A19").Value = -0.012091856045843841
$ws.Range("B19").Value = 0.012017180518879123
$ws.Range("A20").Value = -0.0080171805270818908
$ws.Range("B20").Value = 0.0080057081102946626
$ws.Range("A21").Value = -0.0040057081185711496
$ws.Range("B21").Value = 0.003999999991654235
$ws.Range("A22").Value = -0.045703728164594892
$ws.Range("B22").Value = 0.045492643206083372
$ws.Range("A23").Value = -0.040492643217596047
$ws.Range("B23").Value = 0.040097848822135518
$ws.Range("A24").Value = -0.020097848862280721
$ws.Range("B24").Value = 0.019999999959311232
$ws.Range("A25").Value = -0.0051194097940392425
$ws.Range("B25").Value = 0.0050811829806747255
$ws.Range("A26").Value = -0.0025811829899531347
$ws.Range("B26").Value = 0.0025355996361202671
$ws.Range("A27").Value = -0.000035599645408446179
$ws.Range("B27").Value = -0.00021217883197666154
$ws.Range("A28").Value = 0.0022121788235205386
$ws.Range("B28").Value = -0.0023647659779229002
$ws.Range("A29").Value = -0.0077015333758128079
$ws.Range("B29").Value = 0.0076197877555657101
$ws.Range("A30").Value = 0.052380212129871317
$ws.Range("B30").Value = -0.052707437673345403
$ws.Range("A31").Value = -0.014023480461405669
$ws.Range("B31").Value = 0.014001576280005423
$ws.Range("A32").Value = -0.0040015763014888961
$ws.Range("B32").Value = 0.0039999999894462235
